# Generate Report for Handoff
# Updates the localization-status report: the b.md row moves from
# "Handed back: in sync with en-US" to "Ready for handoff" once a fresh
# handoff xliff has been generated for it (for both zh-cn and de-de),
# and records that the current handback file isn't the latest version.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: update the b.md row's status columns (zh-cn, de-de)
# and the "Latest HO Xliff Generate Date" column.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 08:50:39"

# ---------------------------------------------------------------------
# zh-cn sheet: update the b.md detail row (row 3) with the new handoff
# info, and widen the Error Detail column.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-06 08:50:33"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5bdb6c3b9508e17da294796517176ecc3c505587/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f95276b763906dc35e4ae1e9cc64be0efe3bb2f/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.16

# ---------------------------------------------------------------------
# de-de sheet: same updates as zh-cn, plus widen the Error Detail column.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-06 08:50:39"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5bdb6c3b9508e17da294796517176ecc3c505587/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f95276b763906dc35e4ae1e9cc64be0efe3bb2f/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.16
